$wb = $excel.ActiveWorkbook

# Rename sheets (new participant-generation run ids)
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

$ws1.Name = "GNG_TO-1650291212367313"
$ws2.Name = "NB_TO-16502912148271298"
$ws3.Name = "RS_TO-16502912148280983"
$ws4.Name = "TOL_TO-16502912148751407"
$ws5.Name = "vSAT_TO-1650291214939133"

# Sheet 1 (GNG) - updated stim file names
$ws1.Range("B2").Value = "go_stims-16502912123343117.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912123513103.csv"
$ws1.Range("B4").Value = "go_stims-16502912123533125.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912123663125.csv"

# Sheet 2 (NB) - updated stim file names
$ws2.Range("B2").Value = "ZB-match_4-16502912136001031.csv"
$ws2.Range("B3").Value = "ZB-match_0-16502912135591154.csv"
$ws2.Range("B4").Value = "TB-16502912145400999.csv"
$ws2.Range("B5").Value = "OB-1650291214485101.csv"
$ws2.Range("B6").Value = "TB-16502912148010976.csv"
$ws2.Range("B7").Value = "ZB-match_0-16502912135150976.csv"
$ws2.Range("B8").Value = "OB-16502912144531276.csv"
$ws2.Range("B9").Value = "OB-16502912144241335.csv"
$ws2.Range("B10").Value = "TB-16502912145910983.csv"

# Sheet 3 (RS) - no cell content changes, only sheet name changed above

# Sheet 4 (TOL) - updated stim file names
$ws4.Range("B2").Value = "MM_stims-16502912148431165.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291214830101.csv"
$ws4.Range("B4").Value = "MM_stims-16502912148591006.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912148440998.csv"
$ws4.Range("B6").Value = "MM_stims-16502912148751407.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291214860098.csv"

# Sheet 5 (vSAT) - updated stim file names
$ws5.Range("B2").Value = "vSAT_stims-16502912149071038.csv"
$ws5.Range("B3").Value = "SAT_stims-1650291214879109.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912148910983.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912149231048.csv"
